$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.050.67"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.852.49"
$ws.Range("E3").Value = "  +1.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "707.20"
$ws.Range("E5").Value = "  +1.12%  "

# Row 6 - Solana
Set-TextValue "D6" "172.97"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.850.25"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.44%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.28%  "

# Row 11 - Toncoin
Set-TextValue "D11" "7.32"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.459"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -0.28%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.80"
$ws.Range("E14").Value = "  +0.99%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.499.98"
$ws.Range("E15").Value = "  +1.21%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.845.29"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "71.064.23"

# Row 18 - Polkadot
Set-TextValue "D18" "7.22"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.80%  "

# Row 20 - Chainlink
Set-TextValue "D20" "17.38"
$ws.Range("E20").Value = "  -2.76%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "497.44"
$ws.Range("E21").Value = "  +3.26%  "

# Row 22 - Uniswap
Set-TextValue "D22" "10.68"
$ws.Range("E22").Value = "  -3.76%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +0.39%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  +1.44%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000148"
$ws.Range("E25").Value = "  +2.20%  "

# Row 26 - RenderToken
$ws.Range("E26").Value = "  +1.53%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.23"
$ws.Range("E27").Value = "  -1.54%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "3.22"
$ws.Range("E28").Value = "  +3.68%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  -3.15%  "

# Row 30 - Dai
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.54"
$ws.Range("E31").Value = "  -0.17%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -0.98%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "29.53"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -3.53%  "

# Row 35 - Aptos
Set-TextValue "D35" "9.20"
$ws.Range("E35").Value = "  -0.39%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "3.808.50"

# Row 37 - Binance-PegBSC-USD
$ws.Range("E37").Value = "  +0.23%  "

# Row 38 - Hedera
Set-TextValue "D38" "0.104"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39 - Stacks
Set-TextValue "D39" "2.38"
$ws.Range("E39").Value = "  +7.19%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  +0.30%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  +6.13%  "

# Row 42 - dogwifhat
Set-TextValue "D42" "3.35"
$ws.Range("E42").Value = "  -2.88%  "

# Row 43 - was USDe, now FLOKI
$ws.Range("B43").Value = "FLOKI"
$ws.Range("C43").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D43" "0.000326"
$ws.Range("E43").Value = "  -0.89%  "

# Row 44 - was FirstDigitalUSD, now USDe
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45 - was FLOKI, now FirstDigitalUSD
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D45" "1.00"
$ws.Range("E45").Value = "  +0.10%  "

# Row 46 - Monero
Set-TextValue "D46" "163.30"
$ws.Range("E46").Value = "  +0.58%  "

# Row 47 - OKB
Set-TextValue "D47" "48.67"

# Row 48 - ONDO
$ws.Range("E48").Value = "  +0.36%  "

# Row 49 - Bittensor
Set-TextValue "D49" "416.19"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50 - TheGraph
$ws.Range("E50").Value = "  -1.02%  "

# Row 51 - Cosmos
Set-TextValue "D51" "8.62"
$ws.Range("E51").Value = "  +0.42%  "
